# Update atlas-tactics.xlsx to match atlas-data as of 2024-03-11
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Replace the STIX IDs (column B, rows 2-15) with the refreshed UUIDs
$stixIds = @(
    "x-mitre-tactic--b88656b8-92b5-48c8-aa0d-401d09225c3a",
    "x-mitre-tactic--8407dd06-d298-4fcb-b42c-5e459685d96c",
    "x-mitre-tactic--aea8785b-aff2-4ffb-97a9-724e16802cd0",
    "x-mitre-tactic--38834033-1aaf-4fac-a120-baae566da1f4",
    "x-mitre-tactic--560abc82-7f87-4ae5-8d5a-18e0de351c4e",
    "x-mitre-tactic--86985fad-fc0c-4b01-9441-f1c005dc529e",
    "x-mitre-tactic--f84385a0-14f1-41d0-8256-02e34d3b6fd1",
    "x-mitre-tactic--332194a6-2b0b-445e-900f-28c5588d3996",
    "x-mitre-tactic--eec9cdf5-d82d-4b12-8d7c-5bcf661c4c8c",
    "x-mitre-tactic--7343b749-f8e5-4bc5-88a3-78b80a05456c",
    "x-mitre-tactic--a2b28f78-92b9-476f-968c-427462a3d057",
    "x-mitre-tactic--d7e6e3f1-7033-47a0-b162-51cffda7932a",
    "x-mitre-tactic--0bdff1ba-04fb-43bc-9559-de278508fe94",
    "x-mitre-tactic--d59655a1-d955-47af-bea8-a776ae7383bb"
)

for ($i = 0; $i -lt $stixIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $stixIds[$i]
}

# 2) Fix the typo in the Exfiltration tactic description ("it's" -> "its")
$ws.Range("D7").Value = "The adversary is trying to steal machine learning artifacts or other information about the machine learning system.`n`nExfiltration consists of techniques that adversaries may use to steal data from your network.`nData may be stolen for its valuable intellectual property, or for use in staging future operations.`n`nTechniques for getting data out of a target network typically include transferring it over their command and control channel or an alternate channel and may also include putting size limits on the transmission.`n"

# 3) Update the "created"/"last modified" date string from "31 October 2023" to "11 March 2024"
for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 6).Value = "11 March 2024"
    $ws.Cells.Item($row, 7).Value = "11 March 2024"
}
